$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row => @(old_count, new_count, old_amount, new_amount) -- values applied directly
$updates = @{
    65  = @(2020, 13691790)
    91  = @(151201, 482924639)
    92  = @(409263, 1597003409)
    93  = @(209648, 1309917695)
    94  = @(94226, 918566141)
    96  = @(17317, 796835901)
    98  = @(812, 117970793)
    107 = @(6392, 21970305)
    144 = @(24419, 202157062)
    153 = @(99155, 337748796)
    158 = @(3848, 140904137)
    174 = @(226108, 900747512)
    175 = @(80786, 486199250)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
